$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Effective vs Man" (E) and "Effective vs Zone" (F) to 0.5 for all data rows
$ws.Range("E2:F69").Value = 0.5

# Scroll the frozen pane back to the top (first visible row under the frozen
# header) and reset the selection there instead of the old A70.
$ws.Range("A2").Select()
